$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A5 (was a placeholder "null" value) - row 5 now only has B5
$ws.Range("A5").Clear()

# Row 6 is removed entirely - clear both A6 and B6
$ws.Range("A6:B6").Clear()

# Row 7: A7 now holds "Admin" (previously a placeholder), B7 is cleared
$ws.Range("A7").Value = "Admin"
$ws.Range("B7").Clear()

# Row 8: B8 previously held a placeholder "null" value; now holds a real value
$ws.Range("B8").Value = "dkjhdbsdsk"

# Update sheet view: reset top-left scroll to default, move selection to A7
$ws.Range("A7").Select()
